$d = $word.ActiveDocument

# 1) "From the results above..." paragraph: replace middle clause.
$d.Content.Find.Execute(
    "the net number of collisions in the function",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the ‘clustering’ that occurred in linear probing",
    2) | Out-Null

# 2) "Using a secondary hash function..." paragraph: replace middle clause.
$d.Content.Find.Execute(
    " the values are put in an entirely different part of the map.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " future put attempts are in more distant parts of the map.",
    2) | Out-Null

# 3) "Not all of the keys..." paragraph: rewrite the whole sentence set after
#    "a secondary modulus of 223 ".
$d.Content.Find.Execute(
    "saw a key-entry failure with 1642 hashes stored, while with a secondary modulus of 647 only 1139 hashes were in the array, less than 57% of the maximum capacity of the map. This could be solved by resorting to linear probing if probing with the secondary hash fails, or adding a tertiary hash function with a smaller secondary modulus.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "encountered a key-entry failure after storing 1642 hashes, while a secondary modulus of 647 failed after only 1139 hashes were stored. The reason for this error is that the secondary modulus produced was a factor of the map’s size; adjusting the size of the double hash map to a prime number would negate this issue.",
    2) | Out-Null

# 4) "The results above suggest..." paragraph: rewrite tail.
$d.Content.Find.Execute(
    "a double hash map produces the least conflicts, making it the faster function; however care must be taken to ensure that other methods can be used when the double hashing fails to find a place in the map, or the function can be space-inefficient",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "care must be taken to ensure that when using a double hash map, the map should have a prime-number size to avoid key-entry failures",
    2) | Out-Null

# 5) Move the _GoBack bookmark from its old spot (inside "too fast to
#    impede") to the new spot (inside "discern" -> "discer|n").
$rng = $d.Content
$found = $rng.Find.Execute("discern the password from the hash code", $true)
if ($found) {
    $pos = $rng.Start + 6
    $bmRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
}

# 6) Page margins: top 1076 -> 488 twips, bottom 1701 -> 1083 twips.
$ps = $d.PageSetup
$ps.TopMargin = 24.4
$ps.BottomMargin = 54.15

Write-Output "Done"
